$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.107.76"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.49%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.345.82"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +6.61%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.35%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.68"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +5.19%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.07"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.26%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.643"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +4.06%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.26%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.630"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +7.17%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.75"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.62%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0937"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.02%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.90"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.24%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +13.28%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.106"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.63%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.21"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +9.78%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.699.92"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +6.50%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.345.38"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +5.00%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.084.16"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.87%  "

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +5.10%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.28"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.01%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "75.20"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.46%  "

# Row 22
$ws.Range("B22").Value = "ImmutableX"
$ws.Range("C22").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.57"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +14.88%  "

# Row 23
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.42"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.82%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "252.45"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +11.65%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.03"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.54%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.04"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +5.62%  "

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.12%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.49%  "

# Row 29
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.24"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.61%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.52"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +8.97%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.95"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.81%  "

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.64%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0924"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +6.70%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.95"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +9.51%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.134"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +7.80%  "

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.91%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.76%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.11"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.43%  "

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.11%  "

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +11.16%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.39"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.79%  "

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +16.02%  "

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.97%  "

# Row 44
$ws.Range("B44").Value = "Celestia"
$ws.Range("C44").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.81"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.77%  "

# Row 45
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.09%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.62"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.10%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.37"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +11.76%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "109.46"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +7.84%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.33%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.100"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.97%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "69.93"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +6.49%  "
